$d = $word.ActiveDocument

# The two logo pictures that live in the headers/footers get renamed:
#   - the Pearson Edexcel logo (alt text ends in PearsonLogo.png)   -> image1.png
#   - the BTEC logo            (alt text "BTec_Logo-Orange")        -> image2.jpg
# Look at every section's headers and footers and rename whichever
# inline picture we find using its existing alt text/description to
# decide the new name.

for ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $sec = $d.Sections.Item($secIdx)

    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $shapeCount = $hdr.Range.InlineShapes.Count
            for ($j = 1; $j -le $shapeCount; $j++) {
                $shp = $hdr.Range.InlineShapes.Item($j)
                $desc = $shp.AlternativeText
                $newName = $null
                if ($desc -like "*PearsonLogo.png") {
                    $newName = "image1.png"
                } elseif ($desc -eq "BTec_Logo-Orange") {
                    $newName = "image2.jpg"
                }
                if ($newName) {
                    [void]$shp.Select()
                    $word.Selection.InlineShapes.Item(1).Name = $newName
                }
            }
        }
    }

    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $shapeCount = $ftr.Range.InlineShapes.Count
            for ($j = 1; $j -le $shapeCount; $j++) {
                $shp = $ftr.Range.InlineShapes.Item($j)
                $desc = $shp.AlternativeText
                $newName = $null
                if ($desc -like "*PearsonLogo.png") {
                    $newName = "image1.png"
                } elseif ($desc -eq "BTec_Logo-Orange") {
                    $newName = "image2.jpg"
                }
                if ($newName) {
                    [void]$shp.Select()
                    $word.Selection.InlineShapes.Item(1).Name = $newName
                }
            }
        }
    }
}
